# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.2917716402565462; C = 0.306821227259698;  D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 1.242251378316819 }
    3 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 5.586269137925634 }
    4 = @{ B = 0.2917716402565462; C = 0.306821227259698;  D = 0.1494219747398047; E = 0.4942365360607697; F = 1; G = 1.242251378316819 }
    5 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 5.586269137925634 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
